# The author renamed the sheet "20202s3" to "20202s4" and, while doing so,
# made it the active/selected tab (moving the selection on it to H33),
# which in turn clears the "tabSelected" flag that used to sit on the last
# sheet ("20202s8") -- its own selection (H33) is left untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("20202s3")
$ws.Name = "20202s4"
$ws.Activate()
$ws.Range("H33").Select()
